$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values: B2 + D2 + E2 updated, C2 cleared entirely
$ws.Range("B2").Value = 19.601292356545368
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 19.33729482102439
$ws.Range("E2").Value = 36.578736032412017

# Row 3 values: B3 cleared entirely, C3 + E3 updated, D3 newly added
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 28.126888271470293
$ws.Range("D3").Value = 19.770471625496249
$ws.Range("E3").Value = 24.433638921131998

# Update the active selection to match the new used range of interest
$ws.Range("B1:E3").Select() | Out-Null
